# Actualización desde MV -datos-
# Adds a new "Agosto.2021" column (BH) to the worksheet, copying forward
# the last known value (column BG) for every data row, matching the
# pattern used for every prior "latest data point" column in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label in the new last column (BH == column 60), header style
# (bold/centered/bordered) copied from the previous header cell (BG1).
$ws.Cells.Item(1, 60).Value = "Agosto.2021"
$ws.Cells.Item(1, 59).Copy()
$ws.Cells.Item(1, 60).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Copy the previous last column's (BG == column 59) value forward into
# the new column for every data row (2 through 19).
for ($r = 2; $r -le 19; $r++) {
    $prev = $ws.Cells.Item($r, 59).Value()
    $ws.Cells.Item($r, 60).Value = $prev
}
